$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (shifts H..V to I..W)
$ws.Columns("H").Insert()

# New header cell H1: "DynamiteOccurrenceLinkLocation", bold + text format (matches style used for other headers)
$ws.Range("H1").Value = "DynamiteOccurrenceLinkLocation"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Font.Bold = $true

# New data cell H2: "Main Menu", text format (matches other text data cells in row 2)
$ws.Range("H2").Value = "Main Menu"
$ws.Range("H2").NumberFormat = "@"

# Column width adjustments (closest achievable given engine's width quantization)
$ws.Columns("E").ColumnWidth = 8.333333333333334
$ws.Columns("F").ColumnWidth = 33.333333333333336
$ws.Columns("G").ColumnWidth = 30.333333333333332
$ws.Columns("H").ColumnWidth = 17.666666666666668

# Restore the selection to H5 like in the final workbook
$null = $ws.Range("H5").Select()
